$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price/volume snapshot (Coin, Link, Price, Volume(1h))
# as captured by the scheduled GitHub Actions scraper run.
# Numeric-looking Price values (e.g. "0.999", "314.25") are written with a
# leading apostrophe so Excel keeps them as literal text instead of
# auto-converting them into numbers.

$ws.Range('D2').Value = '41.522.92'
$ws.Range('D3').Value = '2.481.04'
$ws.Range('E3').Value = '  +0.48%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '''314.25'
$ws.Range('E5').Value = '  +0.96%  '
$ws.Range('D6').Value = '''92.40'
$ws.Range('E6').Value = '  -2.33%  '
$ws.Range('D7').Value = '''0.542'
$ws.Range('E7').Value = '  -1.67%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = '''0.501'
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('D10').Value = '''32.53'
$ws.Range('E10').Value = '  -2.89%  '
$ws.Range('D11').Value = '''0.0785'
$ws.Range('E11').Value = '  +0.75%  '
$ws.Range('E12').Value = '  +2.00%  '
$ws.Range('D13').Value = '2.865.72'
$ws.Range('E13').Value = '  +0.49%  '
$ws.Range('D14').Value = '''6.80'
$ws.Range('E14').Value = '  -3.18%  '
$ws.Range('D15').Value = '''15.79'
$ws.Range('E15').Value = '  +7.70%  '
$ws.Range('D16').Value = '2.482.15'
$ws.Range('E16').Value = '  +0.45%  '
$ws.Range('D17').Value = '''0.754'
$ws.Range('E17').Value = '  -3.91%  '
$ws.Range('D18').Value = '41.487.02'
$ws.Range('E18').Value = '  +0.50%  '
$ws.Range('D19').Value = '''6.33'
$ws.Range('E19').Value = '  -0.04%  '
$ws.Range('E20').Value = '  +1.14%  '
$ws.Range('D21').Value = '''71.38'
$ws.Range('E21').Value = '  +4.23%  '
$ws.Range('D22').Value = '''11.16'
$ws.Range('E22').Value = '  -2.02%  '
$ws.Range('D23').Value = '''236.15'
$ws.Range('E23').Value = '  -0.50%  '
$ws.Range('D24').Value = '''2.71'
$ws.Range('E24').Value = '  -2.19%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').Value = '''1.00'
$ws.Range('E25').Value = '  -0.07%  '
$ws.Range('B26').Value = 'ImmutableX'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D26').Value = '''1.89'
$ws.Range('E26').Value = '  -1.51%  '
$ws.Range('D27').Value = '''24.95'
$ws.Range('E27').Value = '  +2.89%  '
$ws.Range('D28').Value = '''2.20'
$ws.Range('E28').Value = '  -0.92%  '
$ws.Range('D29').Value = '''9.65'
$ws.Range('E29').Value = '  -0.11%  '
$ws.Range('D30').Value = '''35.78'
$ws.Range('E30').Value = '  -0.16%  '
$ws.Range('D31').Value = '''157.02'
$ws.Range('E31').Value = '  +3.09%  '
$ws.Range('D32').Value = '''5.43'
$ws.Range('E32').Value = '  -2.27%  '
$ws.Range('D33').Value = '''2.59'
$ws.Range('E33').Value = '  +0.29%  '
$ws.Range('D34').Value = '''0.0749'
$ws.Range('E34').Value = '  +0.81%  '
$ws.Range('D35').Value = '''17.60'
$ws.Range('E35').Value = '  +1.82%  '
$ws.Range('D36').Value = '''2.47'
$ws.Range('E36').Value = '  -6.32%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').Value = '''0.106'
$ws.Range('E37').Value = '  +2.56%  '
$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D38').Value = '''2.91'
$ws.Range('E38').Value = '  -3.39%  '
$ws.Range('D39').Value = '''1.83'
$ws.Range('E39').Value = '  -2.97%  '
$ws.Range('D41').Value = '''4.09'
$ws.Range('E41').Value = '  -4.83%  '
$ws.Range('E42').Value = '  -0.07%  '
$ws.Range('D43').Value = '''19.78'
$ws.Range('E43').Value = '  -4.72%  '
$ws.Range('D44').Value = '1.969.89'
$ws.Range('E44').Value = '  -0.92%  '
$ws.Range('D45').Value = '''0.0283'
$ws.Range('E45').Value = '  -0.40%  '
$ws.Range('D46').Value = '''2.95'
$ws.Range('E46').Value = '  -3.58%  '
$ws.Range('D47').Value = '''8.93'
$ws.Range('E47').Value = '  +2.52%  '
$ws.Range('D48').Value = '2.723.52'
$ws.Range('E48').Value = '  +0.61%  '
$ws.Range('D49').Value = '''97.06'
$ws.Range('E49').Value = '  +0.39%  '
$ws.Range('D50').Value = '''67.73'
$ws.Range('E50').Value = '  -3.31%  '
$ws.Range('D51').Value = '''72.44'
$ws.Range('E51').Value = '  -3.51%  '
